$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1): add P1=14, Q1=15, matching the style of the existing header cells
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2-25: columns I,J,K,L,M,N,O get new values, and add new columns P,Q = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I
    $ws.Cells.Item($r, 10).Value = 2   # J
    $ws.Cells.Item($r, 11).Value = 1   # K
    $ws.Cells.Item($r, 12).Value = 2   # L
    $ws.Cells.Item($r, 13).Value = 2   # M
    $ws.Cells.Item($r, 14).Value = 2   # N
    $ws.Cells.Item($r, 15).Value = 1   # O
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
